$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.093.66"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "3.407.84"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "255.04"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "662.49"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  -6.65%  "
$ws.Range("D8").Value = "0.432"
$ws.Range("E8").Value = "  -6.19%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "3.402.87"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "0.214"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("D13").Value = "42.38"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").Value = "6.46"
$ws.Range("E14").Value = "  +15.78%  "
$ws.Range("D15").Value = "97.839.06"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "4.038.79"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "8.98"
$ws.Range("E18").Value = "  +20.38%  "
$ws.Range("D19").Value = "3.404.99"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").Value = "0.555"
$ws.Range("E20").Value = "  +28.11%  "
$ws.Range("D21").Value = "17.62"
$ws.Range("E21").Value = "  +4.25%  "
$ws.Range("D22").Value = "10.96"
$ws.Range("E22").Value = "  +6.91%  "
$ws.Range("D23").Value = "3.44"
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("D24").Value = "508.63"
$ws.Range("E24").Value = "  -5.50%  "
$ws.Range("D25").Value = "0.0000206"
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("D26").Value = "6.56"
$ws.Range("E26").Value = "  +4.79%  "
$ws.Range("D27").Value = "100.13"
$ws.Range("D28").Value = "12.84"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "3.591.33"
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("D31").Value = "11.53"
$ws.Range("E31").Value = "  +5.30%  "
$ws.Range("E32").Value = "  +5.03%  "
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +15.19%  "
$ws.Range("D36").Value = "0.574"
$ws.Range("E36").Value = "  +7.97%  "
$ws.Range("D37").Value = "29.72"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "7.96"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "1.52"
$ws.Range("E39").Value = "  +14.67%  "
$ws.Range("D40").Value = "533.93"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "0.876"
$ws.Range("E43").Value = "  +7.79%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "9.12"
$ws.Range("E45").Value = "  +18.22%  "
$ws.Range("D46").Value = "5.84"
$ws.Range("E46").Value = "  +17.26%  "
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "1.73"
$ws.Range("E47").Value = "  +17.42%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0426"
$ws.Range("E48").Value = "  +4.17%  "
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").Value = "3.29"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "54.50"
$ws.Range("E51").Value = "  +9.85%  "
